$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Column A's header changes from "when1" to "col_time"; a new column B
# ("col_time_n") is introduced, plus a blank (but bold-styled) column C.
$ws.Range("A1").Value = "col_time"
$ws.Range("B1").Value = "col_time_n"

$hdrFont = $ws.Range("A1:C1").Font
$hdrFont.Bold = $true
$hdrFont.Size = 11

# --- Data rows ----------------------------------------------------------
# Values are Excel serial-day fractions representing times of day:
#   A2/B2 -> 00:01:00 / 00:02:00
#   A3/B3 -> 07:07:07 / 07:07:07
#   A4    -> 07:07:08 (no B4 value)
$ws.Range("A2").Value = 0.0006944444444444445
$ws.Range("B2").Value = 0.001388888888888889
$ws.Range("A3").Value = 0.2966087962995516
$ws.Range("B3").Value = 0.2966087962995516
$ws.Range("A4").Value = 0.29662037037037037

# --- Number formats -------------------------------------------------
# Column A: 12-hour clock with AM/PM; column B: 24-hour h:mm:ss;
# column C: date format applied to the (otherwise empty) cells below
# the header, matching the source sheet's column of blank date cells.
$ws.Range("A2:A4").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"
$ws.Range("B2:B3").NumberFormat = "h:mm:ss"
$ws.Range("C2:C4").NumberFormat = "yyyy\-mm\-dd"

# --- Column widths ----------------------------------------------------
# The emulator's saved <col> width is ColumnWidth + 5/6 (character
# units), so back that constant off to land on the saved widths.
$ws.Columns.Item(1).ColumnWidth = 23.498697916666668
$ws.Columns.Item(2).ColumnWidth = 31.830729166666668

# --- Page setup / selection --------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("A4").Select() | Out-Null
